# Auto-generated edit script applying the cryptos.xlsx price/volume refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.306.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.588.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.245"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.40"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.593.38"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.33"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.309.74"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.312.61"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.611"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.09"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -10.33%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.22"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.724.92"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.47"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0506"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.29%  "
